$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new note about scaling a function (Point struct / x & y scale explanation)
# as a new row below the existing easing-function table.
$ws.Range("A12").Value = "Scaling a function will require you to divide the x values by the y scale, and to multiple the result of a function by the x scale"

# Update the view: select C12 and reset zoom to 100%
[void]$ws.Range("C12").Select()
$excel.ActiveWindow.Zoom = 100
